$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.200.96"
$ws.Range("D3").Value = "'1.637.99"
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'217.04"
$ws.Range("E5").Value = "'  -0.47%  "
$ws.Range("D6").Value = "'0.517"
$ws.Range("E6").Value = "'  +1.91%  "
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E8").Value = "'  -0.11%  "
$ws.Range("E9").Value = "'  +0.20%  "
$ws.Range("D10").Value = "'20.12"
$ws.Range("E10").Value = "'  +0.63%  "
$ws.Range("D11").Value = "'0.0850"
$ws.Range("E11").Value = "'  +0.55%  "
$ws.Range("D12").Value = "'1.866.79"
$ws.Range("E12").Value = "'  -0.12%  "
$ws.Range("D13").Value = "'1.635.85"
$ws.Range("E13").Value = "'  +0.28%  "
$ws.Range("E14").Value = "'  +0.67%  "
$ws.Range("D15").Value = "'0.543"
$ws.Range("E15").Value = "'  +1.71%  "
$ws.Range("E16").Value = "'  -0.86%  "
$ws.Range("D17").Value = "'27.177.27"
$ws.Range("E17").Value = "'  +0.74%  "
$ws.Range("D18").Value = "'0.0₃0739"
$ws.Range("E18").Value = "'  +0.86%  "
$ws.Range("D19").Value = "'217.39"
$ws.Range("E19").Value = "'  -1.12%  "
$ws.Range("E20").Value = "'  +0.08%  "
$ws.Range("D21").Value = "'6.85"
$ws.Range("E21").Value = "'  +1.67%  "
$ws.Range("D22").Value = "'2.53"
$ws.Range("E22").Value = "'  +4.36%  "
$ws.Range("D23").Value = "'4.42"
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("D24").Value = "'9.13"
$ws.Range("E24").Value = "'  -0.59%  "
$ws.Range("D25").Value = "'147.62"
$ws.Range("E25").Value = "'  +0.14%  "
$ws.Range("E26").Value = "'  -0.06%  "
$ws.Range("E27").Value = "'  +0.56%  "
$ws.Range("E28").Value = "'  +0.12%  "
$ws.Range("D29").Value = "'15.67"
$ws.Range("E29").Value = "'  -0.64%  "
$ws.Range("D30").Value = "'0.0507"
$ws.Range("E30").Value = "'  +0.47%  "
$ws.Range("E31").Value = "'  -0.41%  "
$ws.Range("E32").Value = "'  +1.42%  "
$ws.Range("E33").Value = "'  +0.51%  "
$ws.Range("D34").Value = "'1.304.21"
$ws.Range("E34").Value = "'  +2.67%  "
$ws.Range("E35").Value = "'  -0.07%  "
$ws.Range("E36").Value = "'  +0.82%  "
$ws.Range("D37").Value = "'0.0177"
$ws.Range("E37").Value = "'  -0.05%  "
$ws.Range("D38").Value = "'0.549"
$ws.Range("E38").Value = "'  +0.97%  "
$ws.Range("D39").Value = "'0.855"
$ws.Range("E39").Value = "'  +1.67%  "
$ws.Range("E40").Value = "'  -0.01%  "
$ws.Range("D41").Value = "'0.806"
$ws.Range("E41").Value = "'  -0.19%  "
$ws.Range("D42").Value = "'2.25"
$ws.Range("E42").Value = "'  +5.82%  "
$ws.Range("D43").Value = "'5.34"
$ws.Range("E43").Value = "'  -0.37%  "
$ws.Range("D44").Value = "'1.776.26"
$ws.Range("E44").Value = "'  -0.24%  "
$ws.Range("D45").Value = "'62.36"
$ws.Range("E45").Value = "'  +0.18%  "
$ws.Range("D46").Value = "'91.23"
$ws.Range("E46").Value = "'  -1.45%  "
$ws.Range("E47").Value = "'  -1.35%  "
$ws.Range("E48").Value = "'  +2.10%  "
$ws.Range("E49").Value = "'  +0.52%  "
$ws.Range("D50").Value = "'7.65"
$ws.Range("E50").Value = "'  -0.58%  "
$ws.Range("E51").Value = "'  -0.91%  "
